# Scheduled-runner data refresh: updates market-price columns
# (currentAveragePrice/NQ/HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ -> H:N)
# for the affected leve rows across each crafting-class sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 14700
$ws.Cells.Item(21, 9).Value = 5930
$ws.Cells.Item(21, 10).Value = 24444.445
$ws.Cells.Item(21, 11).Value = 5930
$ws.Cells.Item(21, 12).Value = 24444.445
$ws.Cells.Item(21, 13).Value = -5462
$ws.Cells.Item(21, 14).Value = -25380.445

$ws.Cells.Item(23, 8).Value = 14700
$ws.Cells.Item(23, 9).Value = 5930
$ws.Cells.Item(23, 10).Value = 24444.445
$ws.Cells.Item(23, 11).Value = 5930
$ws.Cells.Item(23, 12).Value = 24444.445
$ws.Cells.Item(23, 13).Value = -5696
$ws.Cells.Item(23, 14).Value = -24912.445

$ws.Cells.Item(93, 8).Value = 29978.715
$ws.Cells.Item(93, 10).Value = 29978.715
$ws.Cells.Item(93, 12).Value = 29978.715
$ws.Cells.Item(93, 14).Value = -34970.715

$ws.Cells.Item(106, 8).Value = 41671044
$ws.Cells.Item(106, 9).Value = 2455.8333
$ws.Cells.Item(106, 11).Value = 2455.8333
$ws.Cells.Item(106, 13).Value = -1824.8333

$ws.Cells.Item(111, 8).Value = 1020
$ws.Cells.Item(111, 10).Value = 860
$ws.Cells.Item(111, 12).Value = 2580
$ws.Cells.Item(111, 14).Value = -8714

$ws.Cells.Item(112, 8).Value = 1236.1389
$ws.Cells.Item(112, 10).Value = 1312.1515
$ws.Cells.Item(112, 12).Value = 3936.4545
$ws.Cells.Item(112, 14).Value = -6152.4545

$ws.Cells.Item(132, 8).Value = 33668230
$ws.Cells.Item(132, 9).Value = 47620224
$ws.Cells.Item(132, 10).Value = 1113579
$ws.Cells.Item(132, 11).Value = 142860672
$ws.Cells.Item(132, 12).Value = 3340737
$ws.Cells.Item(132, 13).Value = -142858142
$ws.Cells.Item(132, 14).Value = -3345797

$ws.Cells.Item(137, 8).Value = 2508921.8
$ws.Cells.Item(137, 9).Value = 3402622
$ws.Cells.Item(137, 10).Value = 6560.6
$ws.Cells.Item(137, 11).Value = 10207866
$ws.Cells.Item(137, 12).Value = 19681.8
$ws.Cells.Item(137, 13).Value = -10205316
$ws.Cells.Item(137, 14).Value = -24781.8

$ws.Cells.Item(138, 8).Value = 2821.38
$ws.Cells.Item(138, 9).Value = 1398.8
$ws.Cells.Item(138, 10).Value = 2896.2527
$ws.Cells.Item(138, 11).Value = 4196.4
$ws.Cells.Item(138, 12).Value = 8688.758099999999
$ws.Cells.Item(138, 13).Value = 943.6000000000004
$ws.Cells.Item(138, 14).Value = -18968.7581

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2476.2
$ws.Cells.Item(45, 9).Value = 2190.5
$ws.Cells.Item(45, 11).Value = 2190.5
$ws.Cells.Item(45, 13).Value = -1813.5

$ws.Cells.Item(61, 8).Value = 1764
$ws.Cells.Item(61, 9).Value = 1210.6923
$ws.Cells.Item(61, 11).Value = 1210.6923
$ws.Cells.Item(61, 13).Value = -998.6922999999999

$ws.Cells.Item(74, 8).Value = 5602.091
$ws.Cells.Item(74, 9).Value = 8297.091
$ws.Cells.Item(74, 10).Value = 2907.0908
$ws.Cells.Item(74, 11).Value = 8297.091
$ws.Cells.Item(74, 12).Value = 2907.0908
$ws.Cells.Item(74, 13).Value = -7423.091
$ws.Cells.Item(74, 14).Value = -4655.0908

$ws.Cells.Item(77, 8).Value = 5602.091
$ws.Cells.Item(77, 9).Value = 8297.091
$ws.Cells.Item(77, 10).Value = 2907.0908
$ws.Cells.Item(77, 11).Value = 41485.455
$ws.Cells.Item(77, 12).Value = 14535.454
$ws.Cells.Item(77, 13).Value = -37117.455
$ws.Cells.Item(77, 14).Value = -23271.454

$ws.Cells.Item(102, 8).Value = 3662
$ws.Cells.Item(102, 9).Value = 3470
$ws.Cells.Item(102, 10).Value = 3950
$ws.Cells.Item(102, 11).Value = 3470
$ws.Cells.Item(102, 12).Value = 3950
$ws.Cells.Item(102, 13).Value = -1848
$ws.Cells.Item(102, 14).Value = -7194

$ws.Cells.Item(132, 8).Value = 1775.0303
$ws.Cells.Item(132, 9).Value = 797.94116
$ws.Cells.Item(132, 10).Value = 2813.1875
$ws.Cells.Item(132, 11).Value = 2393.82348
$ws.Cells.Item(132, 12).Value = 8439.5625
$ws.Cells.Item(132, 13).Value = 136.17652
$ws.Cells.Item(132, 14).Value = -13499.5625

$ws.Cells.Item(136, 8).Value = 1764
$ws.Cells.Item(136, 9).Value = 1210.6923
$ws.Cells.Item(136, 11).Value = 3632.0769
$ws.Cells.Item(136, 13).Value = -1082.0769

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1774.9
$ws.Cells.Item(134, 9).Value = 1096.75
$ws.Cells.Item(134, 11).Value = 3290.25
$ws.Cells.Item(134, 13).Value = -755.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5223.2856
$ws.Cells.Item(31, 9).Value = 2056
$ws.Cells.Item(31, 10).Value = 5751.1665
$ws.Cells.Item(31, 11).Value = 2056
$ws.Cells.Item(31, 12).Value = 5751.1665
$ws.Cells.Item(31, 13).Value = -1761
$ws.Cells.Item(31, 14).Value = -6341.1665

$ws.Cells.Item(34, 8).Value = 5223.2856
$ws.Cells.Item(34, 9).Value = 2056
$ws.Cells.Item(34, 10).Value = 5751.1665
$ws.Cells.Item(34, 11).Value = 2056
$ws.Cells.Item(34, 12).Value = 5751.1665
$ws.Cells.Item(34, 13).Value = -1854
$ws.Cells.Item(34, 14).Value = -6155.1665

$ws.Cells.Item(58, 8).Value = 2890.9841
$ws.Cells.Item(58, 9).Value = 1755.3137
$ws.Cells.Item(58, 11).Value = 1755.3137
$ws.Cells.Item(58, 13).Value = -1552.3137

$ws.Cells.Item(99, 8).Value = 12504730
$ws.Cells.Item(99, 9).Value = 33335730
$ws.Cells.Item(99, 10).Value = 6130
$ws.Cells.Item(99, 11).Value = 33335730
$ws.Cells.Item(99, 12).Value = 6130
$ws.Cells.Item(99, 13).Value = -33334232
$ws.Cells.Item(99, 14).Value = -9126

$ws.Cells.Item(126, 8).Value = 12504730
$ws.Cells.Item(126, 9).Value = 33335730
$ws.Cells.Item(126, 10).Value = 6130
$ws.Cells.Item(126, 11).Value = 100007190
$ws.Cells.Item(126, 12).Value = 18390
$ws.Cells.Item(126, 13).Value = -100004720
$ws.Cells.Item(126, 14).Value = -23330

$ws.Cells.Item(132, 8).Value = 4016.087
$ws.Cells.Item(132, 9).Value = 3964.923
$ws.Cells.Item(132, 10).Value = 4082.6
$ws.Cells.Item(132, 11).Value = 11894.769
$ws.Cells.Item(132, 12).Value = 12247.8
$ws.Cells.Item(132, 13).Value = -9364.769
$ws.Cells.Item(132, 14).Value = -17307.8

$ws.Cells.Item(134, 8).Value = 1971.1177
$ws.Cells.Item(134, 9).Value = 1116.1538
$ws.Cells.Item(134, 11).Value = 3348.4614
$ws.Cells.Item(134, 13).Value = -813.4614000000001

$ws.Cells.Item(136, 8).Value = 2890.9841
$ws.Cells.Item(136, 9).Value = 1755.3137
$ws.Cells.Item(136, 11).Value = 5265.9411
$ws.Cells.Item(136, 13).Value = -2715.9411

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 607863.4399999999
$ws.Cells.Item(5, 10).Value = 742805.9
$ws.Cells.Item(5, 12).Value = 2228417.7
$ws.Cells.Item(5, 14).Value = -2228641.7

$ws.Cells.Item(80, 8).Value = 13990.833
$ws.Cells.Item(80, 10).Value = 14231.471
$ws.Cells.Item(80, 12).Value = 42694.413
$ws.Cells.Item(80, 14).Value = -44566.413

$ws.Cells.Item(83, 8).Value = 13990.833
$ws.Cells.Item(83, 10).Value = 14231.471
$ws.Cells.Item(83, 12).Value = 128083.239
$ws.Cells.Item(83, 14).Value = -137443.239

$ws.Cells.Item(107, 8).Value = 67031.266
$ws.Cells.Item(107, 9).Value = 353.33334
$ws.Cells.Item(107, 10).Value = 167048.17
$ws.Cells.Item(107, 11).Value = 1060.00002
$ws.Cells.Item(107, 12).Value = 501144.51
$ws.Cells.Item(107, 13).Value = 859.9999800000001
$ws.Cells.Item(107, 14).Value = -504984.51

$ws.Cells.Item(113, 8).Value = 3906867.2
$ws.Cells.Item(113, 9).Value = 630.2778
$ws.Cells.Item(113, 10).Value = 8929172
$ws.Cells.Item(113, 11).Value = 1890.8334
$ws.Cells.Item(113, 12).Value = 26787516
$ws.Cells.Item(113, 13).Value = 279.1666
$ws.Cells.Item(113, 14).Value = -26791856

$ws.Cells.Item(122, 8).Value = 2506.761
$ws.Cells.Item(122, 9).Value = 765
$ws.Cells.Item(122, 10).Value = 3527.7932
$ws.Cells.Item(122, 11).Value = 6885
$ws.Cells.Item(122, 12).Value = 31750.1388
$ws.Cells.Item(122, 13).Value = -4435
$ws.Cells.Item(122, 14).Value = -36650.1388

$ws.Cells.Item(129, 8).Value = 3896.3
$ws.Cells.Item(129, 10).Value = 3022.1667
$ws.Cells.Item(129, 12).Value = 9066.500100000001
$ws.Cells.Item(129, 14).Value = -19066.5001

$ws.Cells.Item(135, 8).Value = 607863.4399999999
$ws.Cells.Item(135, 10).Value = 742805.9
$ws.Cells.Item(135, 12).Value = 6685253.100000001
$ws.Cells.Item(135, 14).Value = -6690323.100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3003.53
$ws.Cells.Item(126, 9).Value = 2800.1357
$ws.Cells.Item(126, 10).Value = 3870.6316
$ws.Cells.Item(126, 11).Value = 8400.4071
$ws.Cells.Item(126, 12).Value = 11611.8948
$ws.Cells.Item(126, 13).Value = -5930.4071
$ws.Cells.Item(126, 14).Value = -16551.8948

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3656.7917
$ws.Cells.Item(122, 9).Value = 2032.1111
$ws.Cells.Item(122, 11).Value = 6096.3333
$ws.Cells.Item(122, 13).Value = -3646.3333

$ws.Cells.Item(132, 8).Value = 11197.694
$ws.Cells.Item(132, 9).Value = 12683.682
$ws.Cells.Item(132, 10).Value = 8862.571
$ws.Cells.Item(132, 11).Value = 38051.046
$ws.Cells.Item(132, 12).Value = 26587.713
$ws.Cells.Item(132, 13).Value = -35521.046
$ws.Cells.Item(132, 14).Value = -31647.713

$ws.Cells.Item(136, 8).Value = 3775.074
$ws.Cells.Item(136, 9).Value = 925.13336
$ws.Cells.Item(136, 10).Value = 7337.5
$ws.Cells.Item(136, 11).Value = 2775.40008
$ws.Cells.Item(136, 12).Value = 22012.5
$ws.Cells.Item(136, 13).Value = -225.4000800000003
$ws.Cells.Item(136, 14).Value = -27112.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 7938909.5
$ws.Cells.Item(132, 9).Value = 1612.7587
$ws.Cells.Item(132, 10).Value = 25645188
$ws.Cells.Item(132, 11).Value = 4838.2761
$ws.Cells.Item(132, 12).Value = 76935564
$ws.Cells.Item(132, 13).Value = -2308.2761
$ws.Cells.Item(132, 14).Value = -76940624

$ws.Cells.Item(133, 8).Value = 45873
$ws.Cells.Item(133, 10).Value = 45873
$ws.Cells.Item(133, 12).Value = 45873
$ws.Cells.Item(133, 14).Value = -55993

$ws.Cells.Item(136, 8).Value = 5057.9
$ws.Cells.Item(136, 9).Value = 2381.4614
$ws.Cells.Item(136, 10).Value = 10028.429
$ws.Cells.Item(136, 11).Value = 7144.3842
$ws.Cells.Item(136, 12).Value = 30085.287
$ws.Cells.Item(136, 13).Value = -4594.3842
$ws.Cells.Item(136, 14).Value = -35185.287
